$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidate the "Absent" column (H) values to reflect the final report
$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("H14").Value = 0
